$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (shifts old B->C, old C->D)
$ws.Columns("B").Insert()

# New column B gets the same width as column A
$ws.Columns("B").ColumnWidth = $ws.Columns("A").ColumnWidth

# Header for the new column
$ws.Range("B1").Value = "Description"

# Values for the new column
$ws.Range("B2").Value = "LDAP Configuration "
$ws.Range("B3").Value = "Security Token Configuration Check"
$ws.Range("B4").Value = "Value Addco Check"

# Match the wrap-text formatting used by the neighbouring data cells
$ws.Range("B2:B4").WrapText = $true

# Restore the active selection to A2 with no frozen/top-left scrolling
$ws.Range("A2").Select()
